$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new sample row (2026/02/16, 月, 4, 201) was recorded and inserted right
# before the current row 818, pushing the existing rows 818-859 down to
# 819-860 (dimension grows from D859 to D860).
$ws.Rows.Item(818).Insert()

# Leading apostrophe forces the date-like text to stay a literal string
# instead of being auto-parsed into a date serial; ClearFormats() then
# strips the "quote prefix" style marker so the cell ends up with the same
# (default) formatting as its neighbours.
$ws.Cells.Item(818, 1).Value = "'2026/02/16"
$ws.Cells.Item(818, 1).ClearFormats()
$ws.Cells.Item(818, 2).Value = "月"
$ws.Cells.Item(818, 3).Value = 4
$ws.Cells.Item(818, 4).Value = 201
